# Insert a new data row at row 38 (shifting the existing rows 38-47 down to 39-48)
# and populate it with the new "Haba" price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(38).Insert()

$ws.Cells.Item(38, 1).Value = 7
$ws.Cells.Item(38, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(38, 3).Value = "Ñuble"
$ws.Cells.Item(38, 4).Value = 44809
$ws.Cells.Item(38, 5).Value = 16
$ws.Cells.Item(38, 6).Value = 100112026
$ws.Cells.Item(38, 7).Value = "Haba"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 60
$ws.Cells.Item(38, 11).Value = 12000
$ws.Cells.Item(38, 12).Value = 13000
$ws.Cells.Item(38, 13).Value = 12500
$ws.Cells.Item(38, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(38, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(38, 16).Value = 500
$ws.Cells.Item(38, 17).Value = 25
$ws.Cells.Item(38, 18).Value = "Hortaliza"
